# "pre v3.2 e v3.3" - fill in the timing results that were missing for
# attempts 4-7, and tag attempts 4-9 as measured on version v3.2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Attempt 4 (row 8) already had its version set to v3.2; just fill the
# minute/second timing that was missing.
$ws.Range("C8").Value = 22
$ws.Range("E8").Value = 35

# Attempt 5 (row 9)
$ws.Range("C9").Value = 24
$ws.Range("E9").Value = 32
$ws.Range("F9").Value = "v3.2"

# Attempt 6 (row 10)
$ws.Range("C10").Value = 21
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = "v3.2"

# Attempt 7 (row 11)
$ws.Range("C11").Value = 19
$ws.Range("E11").Value = 19
$ws.Range("F11").Value = "v3.2"

# Attempts 8 and 9 (rows 12-13): timing not recorded yet, but already
# tagged as v3.2 runs.
$ws.Range("F12").Value = "v3.2"
$ws.Range("F13").Value = "v3.2"

# Leave the selection where editing stopped.
$ws.Range("F11").Select()
